$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A29").Value = 0.686686
$ws.Range("B29").Value = 0.91091
$ws.Range("C29").Value = 0.6335646007320983
$ws.Range("D29").Value = "query"

$ws.Range("A30").Value = 0.689689
$ws.Range("B30").Value = 0.922922
$ws.Range("C30").Value = 0.6537665147703261
$ws.Range("D30").Value = "query"
